$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell in this sheet (incl. the numeric-looking Price column D) is
# stored as plain text in the source workbook - e.g. "46.50" and "1.000"
# carry a significant trailing zero, and some prices use dotted thousands
# grouping ("22.409.87"). Force Text format on column D before assigning
# so Excel does not auto-detect/convert these strings to numbers (which
# would silently normalize away trailing zeros), keeping cell content and
# type identical to the original text cells.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.409.87'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.566.51'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '286.25'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3712'
$ws.Range('E7').Value = '  +2.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3279'
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.50'
$ws.Range('E9').Value = '  -4.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.146'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07420'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.47'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.846'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.812'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.574.15'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06701'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '86.18'
$ws.Range('E19').Value = '  -2.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.335'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.27'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.78'
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.412.49'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.310'
$ws.Range('E25').Value = '  -2.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.570'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.82'
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.35'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.950'
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.63'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.745.40'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.050'
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.967'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.943'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.679'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08202'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.314'
$ws.Range('E37').Value = '  +2.28%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02383'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06306'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2189'
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.246'
$ws.Range('E41').Value = '  -2.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.12'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6118'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.68'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5950'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.010'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.63'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.180'
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07150'
$ws.Range('E51').Value = '  -0.79%  '
